$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.799.58"
$ws.Range("E2").Value2 = "  +2.15%  "
$ws.Range("D3").Value2 = "1.860.70"
$ws.Range("E3").Value2 = "  +1.67%  "
$ws.Range("D4").Value2 = "'0.9997"
$ws.Range("E4").Value2 = "  +0.02%  "
$ws.Range("D5").Value2 = "'244.95"
$ws.Range("E5").Value2 = "  +0.74%  "
$ws.Range("D6").Value2 = "'0.6403"
$ws.Range("E6").Value2 = "  +3.47%  "
$ws.Range("D7").Value2 = "'1.000"
$ws.Range("E7").Value2 = "  -0.10%  "
$ws.Range("D8").Value2 = "'47.52"
$ws.Range("E8").Value2 = "  +4.14%  "
$ws.Range("D9").Value2 = "'0.07546"
$ws.Range("E9").Value2 = "  +2.71%  "
$ws.Range("D10").Value2 = "'0.2981"
$ws.Range("E10").Value2 = "  +2.49%  "
$ws.Range("E11").Value2 = "  +6.22%  "
$ws.Range("D12").Value2 = "'0.07674"
$ws.Range("E12").Value2 = "  +0.56%  "
$ws.Range("D13").Value2 = "1.856.77"
$ws.Range("E13").Value2 = "  +1.44%  "
$ws.Range("D14").Value2 = "'5.048"
$ws.Range("E14").Value2 = "  +1.53%  "
$ws.Range("D15").Value2 = "'0.6931"
$ws.Range("E15").Value2 = "  +3.39%  "
$ws.Range("D16").Value2 = "'84.01"
$ws.Range("E16").Value2 = "  +2.03%  "
$ws.Range("D17").Value2 = "'0.000009821"
$ws.Range("E17").Value2 = "  +9.39%  "
$ws.Range("D18").Value2 = "'6.130"
$ws.Range("E18").Value2 = "  +4.98%  "
$ws.Range("D19").Value2 = "29.807.63"
$ws.Range("E19").Value2 = "  +2.21%  "
$ws.Range("D20").Value2 = "2.110.57"
$ws.Range("E20").Value2 = "  +1.62%  "
$ws.Range("D21").Value2 = "'236.69"
$ws.Range("E21").Value2 = "  +0.51%  "
$ws.Range("E22").Value2 = "  +1.64%  "
$ws.Range("E23").Value2 = "  -0.11%  "
$ws.Range("D24").Value2 = "'7.506"
$ws.Range("E24").Value2 = "  +2.19%  "
$ws.Range("E25").Value2 = "  -0.03%  "
$ws.Range("D26").Value2 = "'158.95"
$ws.Range("E26").Value2 = "  +0.27%  "
$ws.Range("E27").Value2 = "  +2.39%  "
$ws.Range("D28").Value2 = "'8.568"
$ws.Range("E28").Value2 = "  +0.58%  "
$ws.Range("D29").Value2 = "'17.94"
$ws.Range("E29").Value2 = "  +1.91%  "
$ws.Range("D30").Value2 = "'0.06209"
$ws.Range("E30").Value2 = "  +6.06%  "
$ws.Range("D31").Value2 = "'1.497"
$ws.Range("E31").Value2 = "  +0.46%  "
$ws.Range("D32").Value2 = "'1.300"
$ws.Range("E32").Value2 = "  +5.86%  "
$ws.Range("D33").Value2 = "'4.162"
$ws.Range("E33").Value2 = "  +1.93%  "
$ws.Range("D34").Value2 = "'4.120"
$ws.Range("E34").Value2 = "  +0.94%  "
$ws.Range("D35").Value2 = "'1.901"
$ws.Range("E35").Value2 = "  +2.28%  "
$ws.Range("D36").Value2 = "'1.176"
$ws.Range("E36").Value2 = "  +3.38%  "
$ws.Range("D37").Value2 = "'0.7310"
$ws.Range("E37").Value2 = "  +0.86%  "
$ws.Range("D38").Value2 = "'2.609"
$ws.Range("E38").Value2 = "  -0.02%  "
$ws.Range("D39").Value2 = "'2.824"
$ws.Range("E39").Value2 = "  -1.37%  "
$ws.Range("D40").Value2 = "'0.01788"
$ws.Range("E40").Value2 = "  +1.67%  "
$ws.Range("D41").Value2 = "1.213.02"
$ws.Range("E41").Value2 = "  -1.37%  "
$ws.Range("D42").Value2 = "'0.9232"
$ws.Range("E42").Value2 = "  +2.00%  "
$ws.Range("D43").Value2 = "'6.293"
$ws.Range("E43").Value2 = "  +1.61%  "
$ws.Range("D44").Value2 = "'1.0000"
$ws.Range("E44").Value2 = "  -0.14%  "
$ws.Range("D45").Value2 = "2.020.80"
$ws.Range("E45").Value2 = "  +2.05%  "
$ws.Range("D46").Value2 = "'102.12"
$ws.Range("E46").Value2 = "  +0.27%  "
$ws.Range("D47").Value2 = "'66.80"
$ws.Range("E47").Value2 = "  +1.50%  "
$ws.Range("B48").Value2 = "EnergySwap"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value2 = "'9.228"
$ws.Range("E48").Value2 = "  +0.83%  "
$ws.Range("B49").Value2 = "TheSandbox"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value2 = "'0.4061"
$ws.Range("E49").Value2 = "  +0.42%  "
$ws.Range("B50").Value2 = "RenderToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value2 = "'1.666"
$ws.Range("E50").Value2 = "  +5.18%  "
$ws.Range("B51").Value2 = "Cronos"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value2 = "'0.05795"
$ws.Range("E51").Value2 = "  +0.80%  "
